$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: swap row pairs (full data swap except column A which is the sequential index)
# swap row 142 <-> row 143
$ws.Range("B142").Value2 = 6899181
$ws.Range("B143").Value2 = 6893145
$ws.Range("C142").Value2 = "Portugal Segunda Liga"
$ws.Range("C143").Value2 = "Portugal Segunda Liga"
$ws.Range("D142").Value2 = "Portugal Segunda Liga"
$ws.Range("D143").Value2 = "Portugal Segunda Liga"
$ws.Range("E142").Value2 = 45297.625
$ws.Range("E143").Value2 = 45297.625
$ws.Range("F142").Value2 = "Academico Viseu"
$ws.Range("F143").Value2 = "Nacional"
$ws.Range("G142").Value2 = "UD Leiria"
$ws.Range("G143").Value2 = "Tondela"
$ws.Range("H142").Value2 = 1
$ws.Range("H143").Value2 = 1
$ws.Range("I142").Value2 = 0
$ws.Range("I143").Value2 = 1
$ws.Range("J142").Value2 = "H"
$ws.Range("J143").Value2 = "D"
$ws.Range("K142").Value2 = 2.375
$ws.Range("K143").Value2 = 1.95
$ws.Range("L142").Value2 = 3.2
$ws.Range("L143").Value2 = 3.5
$ws.Range("M142").Value2 = 3
$ws.Range("M143").Value2 = 3.75
$ws.Range("N142").Value2 = 2.375
$ws.Range("N143").Value2 = 1.909
$ws.Range("O142").Value2 = 3
$ws.Range("O143").Value2 = 3.6
$ws.Range("P142").Value2 = 3.1
$ws.Range("P143").Value2 = 3.8
$ws.Range("Q142").Value2 = -0.25
$ws.Range("Q143").Value2 = -0.5
$ws.Range("R142").Value2 = 2.025
$ws.Range("R143").Value2 = 1.825
$ws.Range("S142").Value2 = 1.775
$ws.Range("S143").Value2 = 1.975
$ws.Range("T142").Value2 = 2.25
$ws.Range("T143").Value2 = 2.5
$ws.Range("U142").Value2 = 1.975
$ws.Range("U143").Value2 = 1.875
$ws.Range("V142").Value2 = 1.825
$ws.Range("V143").Value2 = 1.925
$ws.Range("W142").Value2 = 1.375
$ws.Range("W143").Value2 = -1
$ws.Range("X142").Value2 = -1
$ws.Range("X143").Value2 = 2.6
$ws.Range("Y142").Value2 = -1
$ws.Range("Y143").Value2 = -1
$ws.Range("Z142").Value2 = 1.025
$ws.Range("Z143").Value2 = -1
$ws.Range("AA142").Value2 = -1
$ws.Range("AA143").Value2 = 0.9750000000000001
$ws.Range("AB142").Value2 = -1
$ws.Range("AB143").Value2 = -1
$ws.Range("AC142").Value2 = 0.825
$ws.Range("AC143").Value2 = 0.925

# swap row 167 <-> row 168
$ws.Range("B167").Value2 = 6893595
$ws.Range("B168").Value2 = 6893154
$ws.Range("C167").Value2 = "Portugal Segunda Liga"
$ws.Range("C168").Value2 = "Portugal Segunda Liga"
$ws.Range("D167").Value2 = "Portugal Segunda Liga"
$ws.Range("D168").Value2 = "Portugal Segunda Liga"
$ws.Range("E167").Value2 = 45319.520833333336
$ws.Range("E168").Value2 = 45319.520833333336
$ws.Range("F167").Value2 = "FC Porto B"
$ws.Range("F168").Value2 = "Nacional"
$ws.Range("G167").Value2 = "CF Os Belenenses"
$ws.Range("G168").Value2 = "SCU Torreense"
$ws.Range("H167").Value2 = 3
$ws.Range("H168").Value2 = 2
$ws.Range("I167").Value2 = 0
$ws.Range("I168").Value2 = 1
$ws.Range("J167").Value2 = "H"
$ws.Range("J168").Value2 = "H"
$ws.Range("K167").Value2 = 1.5
$ws.Range("K168").Value2 = 2.05
$ws.Range("L167").Value2 = 4
$ws.Range("L168").Value2 = 3.3
$ws.Range("M167").Value2 = 5.5
$ws.Range("M168").Value2 = 3.4
$ws.Range("N167").Value2 = 1.65
$ws.Range("N168").Value2 = 1.833
$ws.Range("O167").Value2 = 3.75
$ws.Range("O168").Value2 = 3.4
$ws.Range("P167").Value2 = 4.5
$ws.Range("P168").Value2 = 4.2
$ws.Range("Q167").Value2 = -0.75
$ws.Range("Q168").Value2 = -0.5
$ws.Range("R167").Value2 = 1.85
$ws.Range("R168").Value2 = 1.85
$ws.Range("S167").Value2 = 1.95
$ws.Range("S168").Value2 = 1.95
$ws.Range("T167").Value2 = 2.75
$ws.Range("T168").Value2 = 2.25
$ws.Range("U167").Value2 = 1.95
$ws.Range("U168").Value2 = 1.8
$ws.Range("V167").Value2 = 1.85
$ws.Range("V168").Value2 = 2
$ws.Range("W167").Value2 = 0.6499999999999999
$ws.Range("W168").Value2 = 0.833
$ws.Range("X167").Value2 = -1
$ws.Range("X168").Value2 = -1
$ws.Range("Y167").Value2 = -1
$ws.Range("Y168").Value2 = -1
$ws.Range("Z167").Value2 = 0.8500000000000001
$ws.Range("Z168").Value2 = 0.8500000000000001
$ws.Range("AA167").Value2 = -1
$ws.Range("AA168").Value2 = -1
$ws.Range("AB167").Value2 = 0.475
$ws.Range("AB168").Value2 = 0.8
$ws.Range("AC167").Value2 = -0.5
$ws.Range("AC168").Value2 = -1

# swap row 201 <-> row 202
$ws.Range("B201").Value2 = 6893167
$ws.Range("B202").Value2 = 6893169
$ws.Range("C201").Value2 = "Portugal Segunda Liga"
$ws.Range("C202").Value2 = "Portugal Segunda Liga"
$ws.Range("D201").Value2 = "Portugal Segunda Liga"
$ws.Range("D202").Value2 = "Portugal Segunda Liga"
$ws.Range("E201").Value2 = 45346.520833333336
$ws.Range("E202").Value2 = 45346.520833333336
$ws.Range("F201").Value2 = "FC Porto B"
$ws.Range("F202").Value2 = "Nacional"
$ws.Range("G201").Value2 = "Leixoes"
$ws.Range("G202").Value2 = "Penafiel"
$ws.Range("H201").Value2 = 0
$ws.Range("H202").Value2 = 3
$ws.Range("I201").Value2 = 2
$ws.Range("I202").Value2 = 2
$ws.Range("J201").Value2 = "A"
$ws.Range("J202").Value2 = "H"
$ws.Range("K201").Value2 = 1.833
$ws.Range("K202").Value2 = 1.75
$ws.Range("L201").Value2 = 3.4
$ws.Range("L202").Value2 = 3.4
$ws.Range("M201").Value2 = 4.333
$ws.Range("M202").Value2 = 4.5
$ws.Range("N201").Value2 = 1.75
$ws.Range("N202").Value2 = 1.666
$ws.Range("O201").Value2 = 3.4
$ws.Range("O202").Value2 = 3.5
$ws.Range("P201").Value2 = 4.5
$ws.Range("P202").Value2 = 5
$ws.Range("Q201").Value2 = -0.5
$ws.Range("Q202").Value2 = -0.75
$ws.Range("R201").Value2 = 1.825
$ws.Range("R202").Value2 = 1.875
$ws.Range("S201").Value2 = 1.975
$ws.Range("S202").Value2 = 1.925
$ws.Range("T201").Value2 = 2.25
$ws.Range("T202").Value2 = 2.25
$ws.Range("U201").Value2 = 1.8
$ws.Range("U202").Value2 = 1.775
$ws.Range("V201").Value2 = 2
$ws.Range("V202").Value2 = 2.025
$ws.Range("W201").Value2 = -1
$ws.Range("W202").Value2 = 0.6659999999999999
$ws.Range("X201").Value2 = -1
$ws.Range("X202").Value2 = -1
$ws.Range("Y201").Value2 = 3.5
$ws.Range("Y202").Value2 = -1
$ws.Range("Z201").Value2 = -1
$ws.Range("Z202").Value2 = 0.4375
$ws.Range("AA201").Value2 = 0.9750000000000001
$ws.Range("AA202").Value2 = -0.5
$ws.Range("AB201").Value2 = -0.5
$ws.Range("AB202").Value2 = 0.7749999999999999
$ws.Range("AC201").Value2 = 0.5
$ws.Range("AC202").Value2 = -1

# swap row 222 <-> row 223
$ws.Range("B222").Value2 = 6893179
$ws.Range("B223").Value2 = 6893607
$ws.Range("C222").Value2 = "Portugal Segunda Liga"
$ws.Range("C223").Value2 = "Portugal Segunda Liga"
$ws.Range("D222").Value2 = "Portugal Segunda Liga"
$ws.Range("D223").Value2 = "Portugal Segunda Liga"
$ws.Range("E222").Value2 = 45360.520833333336
$ws.Range("E223").Value2 = 45360.520833333336
$ws.Range("F222").Value2 = "Benfica B"
$ws.Range("F223").Value2 = "FC Porto B"
$ws.Range("G222").Value2 = "Santa Clara"
$ws.Range("G223").Value2 = "Vilaverdense"
$ws.Range("H222").Value2 = 0
$ws.Range("H223").Value2 = 2
$ws.Range("I222").Value2 = 1
$ws.Range("I223").Value2 = 1
$ws.Range("J222").Value2 = "A"
$ws.Range("J223").Value2 = "H"
$ws.Range("K222").Value2 = 3.2
$ws.Range("K223").Value2 = 1.5
$ws.Range("L222").Value2 = 3.1
$ws.Range("L223").Value2 = 4
$ws.Range("M222").Value2 = 2.2
$ws.Range("M223").Value2 = 5.5
$ws.Range("N222").Value2 = 3.5
$ws.Range("N223").Value2 = 1.571
$ws.Range("O222").Value2 = 3.1
$ws.Range("O223").Value2 = 3.8
$ws.Range("P222").Value2 = 2.1
$ws.Range("P223").Value2 = 5
$ws.Range("Q222").Value2 = 0.25
$ws.Range("Q223").Value2 = -1
$ws.Range("R222").Value2 = 2.025
$ws.Range("R223").Value2 = 1.975
$ws.Range("S222").Value2 = 1.775
$ws.Range("S223").Value2 = 1.825
$ws.Range("T222").Value2 = 2.25
$ws.Range("T223").Value2 = 2.75
$ws.Range("U222").Value2 = 2.025
$ws.Range("U223").Value2 = 2
$ws.Range("V222").Value2 = 1.775
$ws.Range("V223").Value2 = 1.8
$ws.Range("W222").Value2 = -1
$ws.Range("W223").Value2 = 0.571
$ws.Range("X222").Value2 = -1
$ws.Range("X223").Value2 = -1
$ws.Range("Y222").Value2 = 1.1
$ws.Range("Y223").Value2 = -1
$ws.Range("Z222").Value2 = -1
$ws.Range("Z223").Value2 = 0
$ws.Range("AA222").Value2 = 0.7749999999999999
$ws.Range("AA223").Value2 = 0
$ws.Range("AB222").Value2 = -1
$ws.Range("AB223").Value2 = 0.5
$ws.Range("AC222").Value2 = 0.7749999999999999
$ws.Range("AC223").Value2 = -0.5

# swap row 241 <-> row 242
$ws.Range("B241").Value2 = 6893184
$ws.Range("B242").Value2 = 6899271
$ws.Range("C241").Value2 = "Portugal Segunda Liga"
$ws.Range("C242").Value2 = "Portugal Segunda Liga"
$ws.Range("D241").Value2 = "Portugal Segunda Liga"
$ws.Range("D242").Value2 = "Portugal Segunda Liga"
$ws.Range("E241").Value2 = 45381.520833333336
$ws.Range("E242").Value2 = 45381.520833333336
$ws.Range("F241").Value2 = "FC Porto B"
$ws.Range("F242").Value2 = "AVS"
$ws.Range("G241").Value2 = "Penafiel"
$ws.Range("G242").Value2 = "Santa Clara"
$ws.Range("H241").Value2 = 3
$ws.Range("H242").Value2 = 1
$ws.Range("I241").Value2 = 1
$ws.Range("I242").Value2 = 2
$ws.Range("J241").Value2 = "H"
$ws.Range("J242").Value2 = "A"
$ws.Range("K241").Value2 = 1.85
$ws.Range("K242").Value2 = 2.75
$ws.Range("L241").Value2 = 3.75
$ws.Range("L242").Value2 = 3.1
$ws.Range("M241").Value2 = 4
$ws.Range("M242").Value2 = 2.7
$ws.Range("N241").Value2 = 1.833
$ws.Range("N242").Value2 = 3.1
$ws.Range("O241").Value2 = 3.75
$ws.Range("O242").Value2 = 3.1
$ws.Range("P241").Value2 = 4
$ws.Range("P242").Value2 = 2.4
$ws.Range("Q241").Value2 = -0.5
$ws.Range("Q242").Value2 = 0.25
$ws.Range("R241").Value2 = 1.875
$ws.Range("R242").Value2 = 1.75
$ws.Range("S241").Value2 = 1.925
$ws.Range("S242").Value2 = 2.05
$ws.Range("T241").Value2 = 2.5
$ws.Range("T242").Value2 = 2
$ws.Range("U241").Value2 = 1.975
$ws.Range("U242").Value2 = 1.9
$ws.Range("V241").Value2 = 1.825
$ws.Range("V242").Value2 = 1.9
$ws.Range("W241").Value2 = 0.833
$ws.Range("W242").Value2 = -1
$ws.Range("X241").Value2 = -1
$ws.Range("X242").Value2 = -1
$ws.Range("Y241").Value2 = -1
$ws.Range("Y242").Value2 = 1.4
$ws.Range("Z241").Value2 = 0.875
$ws.Range("Z242").Value2 = -1
$ws.Range("AA241").Value2 = -1
$ws.Range("AA242").Value2 = 1.05
$ws.Range("AB241").Value2 = 0.9750000000000001
$ws.Range("AB242").Value2 = 0.8999999999999999
$ws.Range("AC241").Value2 = -1
$ws.Range("AC242").Value2 = -1

# Step 2: update odds on row 245 (no row shift, values refreshed)
$ws.Range("N245").Value2 = 2.1
$ws.Range("P245").Value2 = 2.9
$ws.Range("R245").Value2 = 1.925
$ws.Range("S245").Value2 = 1.925
$ws.Range("U245").Value2 = 1.9
$ws.Range("V245").Value2 = 1.95

# Step 3: insert a new match row at position 246 (shifts old rows 246-252 down to 247-253)
$ws.Rows.Item(246).Insert()
$ws.Range("A245").Copy($ws.Range("A246"))
$ws.Range("E245").Copy($ws.Range("E246"))

# Step 3b: populate new row 246 with the new match data
$ws.Range("A246").Value2 = 244
$ws.Range("B246").Value2 = 6899157
$ws.Range("C246").Value2 = "Portugal Segunda Liga"
$ws.Range("D246").Value2 = "Portugal Segunda Liga"
$ws.Range("E246").Value2 = 45388.29166666666
$ws.Range("F246").Value2 = "Penafiel"
$ws.Range("G246").Value2 = "CF Os Belenenses"
$ws.Range("K246").Value2 = 2
$ws.Range("L246").Value2 = 3.2
$ws.Range("M246").Value2 = 3.3
$ws.Range("N246").Value2 = 1.909
$ws.Range("O246").Value2 = 3.25
$ws.Range("P246").Value2 = 3.6
$ws.Range("Q246").Value2 = -0.5
$ws.Range("R246").Value2 = 1.975
$ws.Range("S246").Value2 = 1.875
$ws.Range("T246").Value2 = 2.25
$ws.Range("U246").Value2 = 1.925
$ws.Range("V246").Value2 = 1.925
$ws.Range("W246").Value2 = 0
$ws.Range("X246").Value2 = 0
$ws.Range("Y246").Value2 = 0
$ws.Range("Z246").Value2 = 0
$ws.Range("AA246").Value2 = 0

# Step 3c: small odds corrections on rows that shifted down (now at their new row numbers)
# row 247 corrections
$ws.Range("N247").Value2 = 2.4
# row 248 corrections
$ws.Range("R248").Value2 = 2.05
$ws.Range("S248").Value2 = 1.8
# row 249 corrections
$ws.Range("N249").Value2 = 2.625
$ws.Range("P249").Value2 = 2.375
$ws.Range("R249").Value2 = 2.05
$ws.Range("S249").Value2 = 1.8
$ws.Range("U249").Value2 = 2
$ws.Range("V249").Value2 = 1.85
# row 250 corrections
$ws.Range("R250").Value2 = 2.025
$ws.Range("S250").Value2 = 1.825
$ws.Range("U250").Value2 = 2.025
$ws.Range("V250").Value2 = 1.825
# row 252 corrections
$ws.Range("R252").Value2 = 1.775
$ws.Range("S252").Value2 = 2.1
$ws.Range("U252").Value2 = 1.9
$ws.Range("V252").Value2 = 1.95